$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.303.22"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.568.01"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.95"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.00"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.567.06"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.81"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.172.61"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.36"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.576.83"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.328.04"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.48"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.63"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.70"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.709.61"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.96"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.46"
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("B33").Value = "RenzoRestakedETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.561.77"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.47"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  -5.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.83"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.63"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.80"
$ws.Range("E40").Value = "  +2.90%  "
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.02"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +3.73%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.24"
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.51"
$ws.Range("E51").Value = "  +3.91%  "
